# Add test row for Treynor_Ratio macro
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "Treynor Ratio1"
$ws.Range("B40").Value = "Test Treynor Ratio"
$ws.Range("C40").Value = "Treynor_Ratio_test1"

$ws.Range("C40").Select()
